# pred_par.xlsx edit: add SVR parameter columns (K:M) and clean up
# redundant cell-level formatting on several cells/rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the three new SVR parameter columns ----------------------------
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0.2
$ws.Range("M2").Value = 20

# --- Normalize cell formatting that reused a redundant duplicate style --
$plainCells = @(
    "B1", "C1", "H1",
    "H2",
    "A5", "A6", "A8", "A11",
    "D14", "E14", "F14", "G14", "I14",
    "D15", "E15", "F15", "G15", "I15",
    "D16", "E16", "F16", "G16", "I16"
)
foreach ($cellAddr in $plainCells) {
    $ws.Range($cellAddr).Style = "Normal"
}

# Row 7 loses its row-level custom formatting
$ws.Rows(7).ClearFormats()

# A13 was an empty, formatted-only cell; remove it completely
$ws.Range("A13").Clear()

# --- Update the selection to match the saved view -----------------------
$ws.Range("K7").Select() | Out-Null
